# process_avatar/final.xlsx deploy edit
#
# The "pic" column (E) on the "S" worksheet stores relative paths to each
# member's avatar image. They used to be rooted from the process_avatar/
# working directory ("../../assets/img/..."), but the generated site page
# only needs the path relative to assets/img/, so the "../../assets/img/"
# prefix is stripped from every picture path in the table.
#
# Besides the data fix, the sheet view is reset: the previous session had
# scrolled down to row 166 with L179 selected; the new save scrolls back to
# the top and leaves K15 selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$prefix = "../../assets/img/"
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 233) { $lastRow = 233 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 5)
    $val = $cell.Value()
    if ($val -ne $null -and $val.ToString().StartsWith($prefix)) {
        $cell.Value = $val.ToString().Substring($prefix.Length)
    }
}

# Reset the view: scroll to top-left and select K15 (matches the new save).
$ws.Range("A1").Select() | Out-Null
$ws.Range("K15").Select() | Out-Null
